# Auto-generated edit script: update crypto price table (columns D and E)
# for rows 2-51, reflecting refreshed market data from the GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "28.012.07"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -1.33%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.765.78"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -3.41%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.65%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "320.77"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.95%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.82%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4251"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -4.33%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3620"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -4.55%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "43.87"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -3.44%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.07456"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -4.03%  "
$ws.Range("E11").Value = "  -3.34%  "
$ws.Range("E12").Value = "  +0.79%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "21.21"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -4.66%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.085"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -3.52%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.334"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -2.77%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.788.65"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.48%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "91.02"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.53%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.00001058"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -2.45%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06405"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.46%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.00%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "17.05"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -3.00%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.979"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -5.74%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "28.017.56"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.40%  "
$ws.Range("E24").Value = "  -3.93%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.139"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.88%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "157.75"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +2.63%  "
$ws.Range("E27").Value = "  -3.84%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.985.03"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.89%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.139"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -10.18%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "124.60"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -4.05%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.165"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -4.15%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "5.652"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -3.71%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.08889"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.98%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.559"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -2.64%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "12.57"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -2.58%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.02320"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -1.41%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.2111"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -3.86%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.06056"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -2.93%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "5.024"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -2.91%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.6378"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -3.68%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.186"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.21%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.9998"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.89%  "
$ws.Range("E43").Value = "  -0.90%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "7.844"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -2.73%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "13.52"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -3.25%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.5920"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -3.24%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.688"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.42%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "122.61"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -3.67%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.188"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.43%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.06868"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -2.09%  "
